$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet: swap SFORZA / TORENBEEK_1982 Xcg comparison values ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
$wsFuselage.Range("C23").Value = 16.834499999999995
$wsFuselage.Range("C24").Value = 17.143322222222217

# --- WING sheet: swap SFORZA / TORENBEEK_1982 Xcg and Ycg comparison values ---
$wsWing = $wb.Worksheets.Item("WING")
$wsWing.Range("C23").Value = 3.5939754358446514
$wsWing.Range("C24").Value = 4.3631082000119275
$wsWing.Range("C27").Value = 6.114221148470394
$wsWing.Range("C28").Value = 4.998846772296348
